$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price / volume(1h) table with the latest scraped values.
#
# Rows 43 and 44 swap coin identities (VeChain <-> Maker) along with updated
# price / volume figures. Row 43 becomes Maker, row 44 becomes VeChain.
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.458.73"
$ws.Range("E43").Value = "  -0.96%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0213"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.34%  "

$ws.Range("D2").Value = '37.689.65'
$ws.Range("E2").Value = '  +1.08%  '
$ws.Range("D3").Value = '2.091.88'
$ws.Range("E3").Value = '  +1.57%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.623'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.13%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '57.71'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.42%  '
$ws.Range("E9").Value = '  +1.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0778'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.34%  '
$ws.Range("E11").Value = '  +2.76%  '
$ws.Range("D12").Value = '2.388.01'
$ws.Range("E12").Value = '  +1.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.43'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.98%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.05'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.765'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.39%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.22'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.75%  '
$ws.Range("D17").Value = '2.077.87'
$ws.Range("E17").Value = '  +0.83%  '
$ws.Range("D18").Value = '37.661.69'
$ws.Range("E18").Value = '  +1.26%  '
$ws.Range("E19").Value = '  -2.74%  '
$ws.Range("E20").Value = '  +1.75%  '
$ws.Range("E21").Value = '  +1.50%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '228.02'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.88%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.41'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.18%  '
$ws.Range("E25").Value = '  -0.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.75'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.87%  '
$ws.Range("E27").Value = '  +9.84%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.93'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.98%  '
$ws.Range("E29").Value = '  -1.51%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.43'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.08%  '
$ws.Range("E31").Value = '  +1.13%  '
$ws.Range("E32").Value = '  +4.06%  '
$ws.Range("E33").Value = '  +1.32%  '
$ws.Range("E34").Value = '  -0.32%  '
$ws.Range("E35").Value = '  +0.20%  '
$ws.Range("E36").Value = '  +4.57%  '
$ws.Range("E37").Value = '  +5.20%  '
$ws.Range("E38").Value = '  +0.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.39'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.85%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0992'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.29%  '
$ws.Range("E41").Value = '  -0.44%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '97.93'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.02%  '
$ws.Range("E45").Value = '  -0.82%  '
$ws.Range("E46").Value = '  +3.66%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '15.57'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.37%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.03'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.71%  '
$ws.Range("E49").Value = '  +2.85%  '
$ws.Range("E50").Value = '  +1.97%  '
$ws.Range("D51").Value = '2.284.65'
$ws.Range("E51").Value = '  +1.53%  '
